$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-22, replacing the old Strike# (pitch count)
# derived values with the new K (strikeout) based values.
$kValues = @{
    2  = 6
    3  = 2
    4  = 5
    5  = 3
    6  = 11
    7  = 7
    8  = 8
    9  = 4
    10 = 3
    11 = 5
    12 = 3
    13 = 1
    14 = 10
    15 = 1
    16 = 6
    17 = 5
    18 = 4
    19 = 7
    20 = 4
    21 = 5
    22 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
